$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 data: ground banana terminal block
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "MKDSN"
$ws.Range("C11").Value = "MKDSN"
$ws.Range("D11").Value = "MKDSN"
$ws.Range("E11").Value = "X1"
$ws.Range("F11").Value = "Terminal Block"
$ws.Range("G11").Value = "PHOENIX"
$ws.Range("H11").Value = 1729128

# H column alignment - set left horizontal alignment only on cells that have content
$ws.Range("H1").HorizontalAlignment = -4131
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("H4").HorizontalAlignment = -4131
$ws.Range("H7").HorizontalAlignment = -4131
$ws.Range("H8").HorizontalAlignment = -4131
$ws.Range("H9").HorizontalAlignment = -4131
$ws.Range("H10").HorizontalAlignment = -4131
$ws.Range("H11").HorizontalAlignment = -4131

# Recalculate best-fit column width for column E after adding new data
$ws.Range("E1:E11").Columns.AutoFit() | Out-Null

# Update selection
$ws.Range("E16").Select()
